$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.719288
$ws.Range("H2").Value = 11.157864
$ws.Range("I2").Value = 0.04235839908674209
$ws.Range("J2").Value = 0.04235839908674209
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1278363333333334
$ws.Range("N2").Value = 0.383509
$ws.Range("O2").Value = 0.002480915078704262
$ws.Range("P2").Value = 0.002480915078704262
$ws.Range("Q2").Value = 0.4754601405306668
$ws.Range("R2").Value = 4.279141264776
$ws.Range("S2").Value = 0.0001050875910040713
$ws.Range("T2").Value = 0.0001050875910040713

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.719288
$ws.Range("H3").Value = 11.157864
$ws.Range("I3").Value = 0.04235839908674209
$ws.Range("J3").Value = 0.04235839908674209
$ws.Range("O3").Value = 0.002269935507489869
$ws.Range("P3").Value = 0.002269935507489869
$ws.Range("Q3").Value = 0.4350265209199999
$ws.Range("R3").Value = 3.91523868828
$ws.Range("S3").Value = 0.00009615083412742229
$ws.Range("T3").Value = 0.00009615083412742229

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.719288
$ws.Range("H4").Value = 11.157864
$ws.Range("I4").Value = 0.04235839908674209
$ws.Range("J4").Value = 0.04235839908674209
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04120633333333334
$ws.Range("N4").Value = 0.123619
$ws.Range("O4").Value = 0.0007996898146180199
$ws.Range("P4").Value = 0.0007996898146180199
$ws.Range("Q4").Value = 0.1532582210906667
$ws.Range("R4").Value = 1.379323989816
$ws.Range("S4").Value = 0.00003387358031319288
$ws.Range("T4").Value = 0.00003387358031319288

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.719288
$ws.Range("H5").Value = 11.157864
$ws.Range("I5").Value = 0.04235839908674209
$ws.Range("J5").Value = 0.04235839908674209
$ws.Range("M5").Value = 51.241888
$ws.Range("N5").Value = 153.725664
$ws.Range("O5").Value = 0.9944494595991877
$ws.Range("P5").Value = 0.9944494595991878
$ws.Range("Q5").Value = 190.583339135744
$ws.Range("R5").Value = 1715.250052221696
$ws.Range("S5").Value = 0.0421232870812974
$ws.Range("T5").Value = 0.0421232870812974

# Row 6
$ws.Range("I6").Value = 0.2979256989470644
$ws.Range("J6").Value = 0.2979256989470644
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1278363333333334
$ws.Range("N6").Value = 0.383509
$ws.Range("O6").Value = 0.002480915078704262
$ws.Range("P6").Value = 0.002480915078704262
$ws.Range("Q6").Value = 3.344125315005223
$ws.Range("R6").Value = 30.097127835047
$ws.Range("S6").Value = 0.0007391283588512787
$ws.Range("T6").Value = 0.0007391283588512787

# Row 7
$ws.Range("I7").Value = 0.2979256989470644
$ws.Range("J7").Value = 0.2979256989470644
$ws.Range("O7").Value = 0.002269935507489869
$ws.Range("P7").Value = 0.002269935507489869
$ws.Range("S7").Value = 0.0006762721226336784
$ws.Range("T7").Value = 0.0006762721226336784

# Row 8
$ws.Range("I8").Value = 0.2979256989470644
$ws.Range("J8").Value = 0.2979256989470644
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.04120633333333334
$ws.Range("N8").Value = 0.123619
$ws.Range("O8").Value = 0.0007996898146180199
$ws.Range("P8").Value = 0.0007996898146180199
$ws.Range("Q8").Value = 1.077934096241889
$ws.Range("R8").Value = 9.701406866177
$ws.Range("S8").Value = 0.0002382481469609219
$ws.Range("T8").Value = 0.0002382481469609219

# Row 9
$ws.Range("I9").Value = 0.2979256989470644
$ws.Range("J9").Value = 0.2979256989470644
$ws.Range("M9").Value = 51.241888
$ws.Range("N9").Value = 153.725664
$ws.Range("O9").Value = 0.9944494595991877
$ws.Range("P9").Value = 0.9944494595991878
$ws.Range("Q9").Value = 1340.458462639435
$ws.Range("R9").Value = 12064.12616375491
$ws.Range("S9").Value = 0.2962720503186185
$ws.Range("T9").Value = 0.2962720503186185

# Row 10
$ws.Range("G10").Value = 6.299630666666666
$ws.Range("H10").Value = 18.898892
$ws.Range("I10").Value = 0.07174552491706633
$ws.Range("J10").Value = 0.07174552491706633
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1278363333333334
$ws.Range("N10").Value = 0.383509
$ws.Range("O10").Value = 0.002480915078704262
$ws.Range("P10").Value = 0.002480915078704262
$ws.Range("Q10").Value = 0.805321685780889
$ws.Range("R10").Value = 7.247895172028
$ws.Range("S10").Value = 0.0001779945545963022
$ws.Range("T10").Value = 0.0001779945545963022

# Row 11
$ws.Range("G11").Value = 6.299630666666666
$ws.Range("H11").Value = 18.898892
$ws.Range("I11").Value = 0.07174552491706633
$ws.Range("J11").Value = 0.07174552491706633
$ws.Range("O11").Value = 0.002269935507489869
$ws.Range("P11").Value = 0.002269935507489869
$ws.Range("Q11").Value = 0.7368363009266665
$ws.Range("R11").Value = 6.631526708339999
$ws.Range("S11").Value = 0.000162857714512748
$ws.Range("T11").Value = 0.000162857714512748

# Row 12
$ws.Range("G12").Value = 6.299630666666666
$ws.Range("H12").Value = 18.898892
$ws.Range("I12").Value = 0.07174552491706633
$ws.Range("J12").Value = 0.07174552491706633
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.04120633333333334
$ws.Range("N12").Value = 0.123619
$ws.Range("O12").Value = 0.0007996898146180199
$ws.Range("P12").Value = 0.0007996898146180199
$ws.Range("Q12").Value = 0.2595846811275556
$ws.Range("R12").Value = 2.336262130148
$ws.Range("S12").Value = 0.0000573741655206013
$ws.Range("T12").Value = 0.0000573741655206013

# Row 13
$ws.Range("G13").Value = 6.299630666666666
$ws.Range("H13").Value = 18.898892
$ws.Range("I13").Value = 0.07174552491706633
$ws.Range("J13").Value = 0.07174552491706633
$ws.Range("M13").Value = 51.241888
$ws.Range("N13").Value = 153.725664
$ws.Range("O13").Value = 0.9944494595991877
$ws.Range("P13").Value = 0.9944494595991878
$ws.Range("Q13").Value = 322.8049690626986
$ws.Range("R13").Value = 2905.244721564288
$ws.Range("S13").Value = 0.07134729848243668
$ws.Range("T13").Value = 0.07134729848243668

# Row 14
$ws.Range("G14").Value = 51.62686066666667
$ws.Range("H14").Value = 154.880582
$ws.Range("I14").Value = 0.5879703770491272
$ws.Range("J14").Value = 0.5879703770491272
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.1278363333333334
$ws.Range("N14").Value = 0.383509
$ws.Range("O14").Value = 0.002480915078704262
$ws.Range("P14").Value = 0.002480915078704262
$ws.Range("Q14").Value = 6.599788569137557
$ws.Range("R14").Value = 59.39809712223801
$ws.Range("S14").Value = 0.00145870457425261
$ws.Range("T14").Value = 0.00145870457425261

# Row 15
$ws.Range("G15").Value = 51.62686066666667
$ws.Range("H15").Value = 154.880582
$ws.Range("I15").Value = 0.5879703770491272
$ws.Range("J15").Value = 0.5879703770491272
$ws.Range("O15").Value = 0.002269935507489869
$ws.Range("P15").Value = 0.002269935507489869
$ws.Range("Q15").Value = 6.038535757876666
$ws.Range("R15").Value = 54.34682182088999
$ws.Range("S15").Value = 0.00133465483621602
$ws.Range("T15").Value = 0.00133465483621602

# Row 16
$ws.Range("G16").Value = 51.62686066666667
$ws.Range("H16").Value = 154.880582
$ws.Range("I16").Value = 0.5879703770491272
$ws.Range("J16").Value = 0.5879703770491272
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.04120633333333334
$ws.Range("N16").Value = 0.123619
$ws.Range("O16").Value = 0.0007996898146180199
$ws.Range("P16").Value = 0.0007996898146180199
$ws.Range("Q16").Value = 2.127353629584222
$ws.Range("R16").Value = 19.146182666258
$ws.Range("S16").Value = 0.0004701939218233038
$ws.Range("T16").Value = 0.0004701939218233038

# Row 17
$ws.Range("G17").Value = 51.62686066666667
$ws.Range("H17").Value = 154.880582
$ws.Range("I17").Value = 0.5879703770491272
$ws.Range("J17").Value = 0.5879703770491272
$ws.Range("M17").Value = 51.241888
$ws.Range("N17").Value = 153.725664
$ws.Range("O17").Value = 0.9944494595991877
$ws.Range("P17").Value = 0.9944494595991878
$ws.Range("Q17").Value = 2645.457812072938
$ws.Range("R17").Value = 23809.12030865645
$ws.Range("S17").Value = 0.5847068237168352
$ws.Range("T17").Value = 0.5847068237168354

